$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formats from row 2's specially-styled cells (date format in A,
# percentage format in S) down to row 3 before assigning the new values, so
# that the existing styles (s="1" / s="2") are reused instead of creating
# new style entries.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

$ws.Range("S2").Copy()
$ws.Range("S3").PasteSpecial(-4122)

$ws.Range("A3").Value = 42632.880844907406
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = "Buy"
$ws.Range("D3").Value = 16
$ws.Range("E3").Value = 12103
$ws.Range("F3").Value = 610
$ws.Range("G3").Value = 61
$ws.Range("H3").Value = 37
$ws.Range("I3").Value = 71
$ws.Range("J3").Value = 28
$ws.Range("K3").Value = 16333
$ws.Range("L3").Value = 165
$ws.Range("M3").Value = 102
$ws.Range("N3").Value = 10
$ws.Range("O3").Value = 4
$ws.Range("P3").Value = "Named"
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0.56
$ws.Range("S3").Value = 0.1015
$ws.Range("T3").Value = -0.93
$ws.Range("U3").Value = 2.32
$ws.Range("V3").Value = "N/A"
$ws.Range("W3").Value = 0
